$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out everything from row 11 downward first (rows 11-24), then rewrite the final layout.
$ws.Range("A11:C24").ClearContents()

# --- Row 10: Objetivos/Docentes text cells now hold the professor entry ---
$ws.Range("B10").Value = '5963230 - Leandro Gonçalves de Aguiar'
$ws.Range("C10").Value = '5963230 - Leandro Gonçalves de Aguiar'

# --- Row 11 ---
$ws.Range("A11").Value = 'Objectives:'

# --- Row 12 ---
$ws.Range("A12").Value = 'Docentes responsáveis:'

# --- Row 13 ---
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

# --- Row 14 ---
$ws.Range("A14").Value = 'Short syllabus:'

# --- Row 15 ---
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2012'
$ws.Range("C15").Value = '01/01/2012'

# --- Row 16 ---
$ws.Range("A16").Value = 'Syllabus:'

# --- Row 17 ---
$ws.Range("A17").Value = 'Avaliação:'

# --- Row 18 ---
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5963230 - Leandro Gonçalves de Aguiar'
$ws.Range("C18").Value = '5963230 - Leandro Gonçalves de Aguiar'

# --- Row 19 ---
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Duas provas escritas e eventual apresentação de trabalho.'
$ws.Range("C19").Value = 'Duas provas escritas e eventual apresentação de trabalho.'

# --- Row 20 ---
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'Nota(N) = 50% Prova P1 + 50% Prova P2. Os pesos poderão ser redefinidos caso seja incorporada nota de trabalho.'
$ws.Range("C20").Value = 'Nota(N) = 50% Prova P1 + 50% Prova P2. Os pesos poderão ser redefinidos caso seja incorporada nota de trabalho.'

# --- Row 21 ---
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Média Final = (N + Prova Recuperação)/2'
$ws.Range("C21").Value = 'Média Final = (N + Prova Recuperação)/2'

# --- Row 22 ---
$ws.Range("A22").Value = 'Requisitos:'

# --- Row 23 ---
$ws.Range("B23").Value = 'LOQ4003 -  Cinética Química Aplicada  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOQ4003 -  Cinética Química Aplicada  (Requisito fraco)
'

# --- Row heights ---
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30

# Remove the now-unused trailing row so the sheet dimension shrinks to A1:C23
$ws.Rows.Item(24).Delete()